$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The figures in this sheet are stored as text (not numbers), so force the
# target cells to keep a text format before writing the updated values back.
$ws.Range("B11:D11").NumberFormat = "@"
$ws.Range("B12:D12").NumberFormat = "@"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"

# Row 11: Enterprises density (per 1000 people)
$ws.Range("B11").Value = "30.43"
$ws.Range("C11").Value = "1.37"
$ws.Range("D11").Value = "31.81"

# Row 12: Employment (% of total)
$ws.Range("B12").Value = "45.68"
$ws.Range("C12").Value = "29.43"
$ws.Range("D12").Value = "75.11"

# Row 14: Enterprises (% of total)
$ws.Range("B14").Value = "95.53"
$ws.Range("D14").Value = "99.83"
